$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
Write-Host $ws.Name
